# Update 8th March 2025 2322 Hours
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 24
$ws.Range("B24").Value = 8
$ws.Range("C24").Value = "Rotate Image"
$ws.Range("E24").Value = "LeetCode"

# Row 25
$ws.Range("C25").Value = "Flipping an Image"
$ws.Range("E25").Value = "Bosscoder Academy"

# Row 26
$ws.Range("C26").Value = "Set Matrix Zero"
$ws.Range("E26").Value = "LeetCode"

# Update the selected cell to reflect where the user left off (E27), matching source file
$ws.Range("E27").Select()
